$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Anticipo Ft. Estero B.Cred.T. (IT74*680)"
$ws.Range("B3").Value = "Anticipo Ft. Italia B.Pop.Soft. (IT15*456)"

$ws.Columns.Item(2).ColumnWidth = 41.386666666666667

$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
